$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)  # column 58 = BF
    if ($cell.Text -eq "1-1-2007-08") {
        $cell.Value = "2008-01-01"
    }
}
